$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the team record columns (AD/AE/AF).
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold/centered/bordered style already used by the other
# header cells (e.g. AC1) by copying its format onto the new headers.
$ws.Range("AC1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# Every player row (2-61) gets the same team record: 74 wins, 88 losses, 0 ties.
for ($row = 2; $row -le 61; $row++) {
    $ws.Cells.Item($row, 30).Value = 74
    $ws.Cells.Item($row, 31).Value = 88
    $ws.Cells.Item($row, 32).Value = 0
}
